$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scanner log row (row 2) was re-logged for the "Immuno&hema" subject
# instead of "Internal Medicine 1". Update the Subject cell (B2) value and
# give it the highlighted/centered formatting that the updated log uses.
$rng = $ws.Range("B2")
$rng.Value = "Immuno&hema"

# Highlight style: light-gray fill, centered alignment, 11pt Calibri font.
$rng.Font.Size = 11
$rng.Interior.Color = 15790320
$rng.Interior.PatternColor = 15790320
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

# Reflect that B2 is the active/selected cell, as in the saved workbook.
$rng.Select()
